# feat: add 2022-Q4 data
#
# 1. "总计" (summary) sheet: insert a new row for 2022-Q4 at the top of the
#    data (row 2), pushing the existing quarters down by one row.
# 2. Insert a brand-new worksheet named "2022-Q4" right after "总计" holding
#    the per-fund holding detail for the new quarter (same shape as the
#    other quarterly sheets). It is built by duplicating the "2022-Q3"
#    sheet (so it inherits the right sheetPr/styles/column layout) and then
#    overwriting its contents.

function Set-TextValue {
    # Write $val into ($row,$col) forcing text storage (Excel normally
    # auto-coerces numeric-looking strings like "59.06" into numbers).
    param($ws, $row, $col, $val)
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet — shift existing rows 2..6 down to 3..7 (copy bottom-up
#    so we never overwrite a row before it has been copied), then write
#    the new 2022-Q4 row into row 2.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A6:D6").Copy($total.Range("A7:D7"))
$total.Range("A5:D5").Copy($total.Range("A6:D6"))
$total.Range("A4:D4").Copy($total.Range("A5:D5"))
$total.Range("A3:D3").Copy($total.Range("A4:D4"))
$total.Range("A2:D2").Copy($total.Range("A3:D3"))

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 9
$total.Range("D2").Value = 1.89

# ---------------------------------------------------------------------
# 2) Duplicate "2022-Q3" (placed right before it -> becomes the 2nd tab,
#    right after "总计"), rename it to "2022-Q4", extend the index-column
#    styling to the extra rows we need, then overwrite all the values.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The template only had 4 data rows (rows 2-5); we need 9 (rows 2-10), so
# replicate the styled index-column cell down to rows 6-10.
$q4.Range("A5").Copy($q4.Range("A6:A10"))

$q4.Cells.Item(1, 2).Value = "基金代码"
$q4.Cells.Item(1, 3).Value = "基金名称"
$q4.Cells.Item(1, 4).Value = "基金规模"
$q4.Cells.Item(1, 5).Value = "股票总仓位"
$q4.Cells.Item(1, 6).Value = "仓位占比"
$q4.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4.Cells.Item(1, 8).Value = "仓位排名"

$data = @(
    @(0, "100032", "富国中证红利指数增强A", "59.06", "91.28", "1.94", "1.1458", 4),
    @(1, "013422", "太平智行三个月定期开放混合", "5.99", "87.78", "5.66", "0.3390", 3),
    @(2, "009794", "太平智选一年定期开放股票", "4.64", "88.49", "5.29", "0.2455", 3),
    @(3, "008682", "富国中证红利指数增强C", "5.45", "91.28", "1.94", "0.1057", 4),
    @(4, "512040", "富国中证价值ETF", "3.39", "99.29", "1.08", "0.0366", 4),
    @(5, "016053", "泰康先进材料股票A", "0.16", "89.64", "3.06", "0.0049", 9),
    @(6, "562520", "华夏中证智选1000成长创新策略ETF", "0.38", "96.24", "1.08", "0.0041", 1),
    @(7, "016054", "泰康先进材料股票C", "0.13", "89.64", "3.06", "0.0040", 9),
    @(8, "005770", "信澳中证沪港深高股息精选指数", "0.13", "23.47", "0.66", "0.0009", 4)
)

$r = 2
foreach ($row in $data) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    Set-TextValue $q4 $r 2 $row[1]
    Set-TextValue $q4 $r 3 $row[2]
    Set-TextValue $q4 $r 4 $row[3]
    Set-TextValue $q4 $r 5 $row[4]
    Set-TextValue $q4 $r 6 $row[5]
    Set-TextValue $q4 $r 7 $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# Restore "总计" as the active sheet (it was the original active tab).
$total.Activate()
